$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name (date rolled forward one day)
$ws.Name = "Through 2021-12-27"

# Update the December label in column A, row 13
$ws.Range("A13").Value = "December (through 12-27)"

# Update December row (row 13) values
$ws.Range("B13").Value = 40
$ws.Range("C13").Value = 87
$ws.Range("D13").Value = 106
$ws.Range("E13").Value = 64
$ws.Range("G13").Value = 126
$ws.Range("H13").Value = 166

# Update Total row (row 14) values
$ws.Range("B14").Value = 331
$ws.Range("C14").Value = 650
$ws.Range("D14").Value = 927
$ws.Range("E14").Value = 746
$ws.Range("G14").Value = 1390
$ws.Range("H14").Value = 1809
